$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH130"
$ws.Range("C2").Value = "SOUTHERN AFRICA- THE IMPOSSIBLE NEUTRALITY"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$row2 = $ws.Range("A2:H2")
$row2.Font.Name = "Calibri"
$row2.Font.Size = 10
